# Apply updated crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "64.096.53"
$ws.Cells.Item(2, 5).Value = "  -0.12%  "

$ws.Cells.Item(3, 4).Value = "2.757.54"
$ws.Cells.Item(3, 5).Value = "  -0.91%  "

$ws.Cells.Item(4, 5).Value = "  +0.06%  "

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "574.74"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -2.26%  "

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "159.33"
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -1.34%  "

$ws.Cells.Item(7, 5).Value = "  +0.06%  "

$ws.Cells.Item(8, 5).Value = "  -3.11%  "

$ws.Cells.Item(9, 5).Value = "  -4.04%  "

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.87"
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -13.67%  "

$ws.Cells.Item(11, 5).Value = "  +3.16%  "

$ws.Cells.Item(13, 4).Value = "3.246.99"
$ws.Cells.Item(13, 5).Value = "  -0.87%  "

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "26.99"
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -2.41%  "

$ws.Cells.Item(15, 4).Value = "63.724.85"
$ws.Cells.Item(15, 5).Value = "  -0.57%  "

$ws.Cells.Item(16, 5).Value = "  -5.85%  "

$ws.Cells.Item(17, 4).Value = "2.765.57"
$ws.Cells.Item(17, 5).Value = "  -0.65%  "

$ws.Cells.Item(18, 5).Value = "  -2.41%  "

$ws.Cells.Item(19, 5).Value = "  -5.08%  "

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "359.63"
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -2.22%  "

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.64"
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -6.22%  "

$ws.Cells.Item(22, 5).Value = "  -0.51%  "

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.528"
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -8.08%  "

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "65.09"
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -3.95%  "

$ws.Cells.Item(25, 5).Value = "  -3.85%  "

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.54"
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -4.37%  "

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.02%  "

$ws.Cells.Item(28, 4).Value = "0.0₃0904"
$ws.Cells.Item(28, 5).Value = "  -7.08%  "

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.37"
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +1.08%  "

$ws.Cells.Item(30, 2).Value = "Fetch.AI"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.37"
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +7.36%  "

$ws.Cells.Item(31, 2).Value = "PancakeSwap"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.96"
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -4.35%  "

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "169.92"
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -1.24%  "

$ws.Cells.Item(33, 2).Value = "EthereumClassic"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "20.22"
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -3.19%  "

$ws.Cells.Item(34, 2).Value = "NEARProtocol"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.94"
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -4.81%  "

$ws.Cells.Item(35, 5).Value = "  -2.08%  "

$ws.Cells.Item(36, 5).Value = "  +0.05%  "

$ws.Cells.Item(37, 5).Value = "  -2.11%  "

$ws.Cells.Item(38, 5).Value = "  -2.35%  "

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "349.50"
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +1.90%  "

$ws.Cells.Item(40, 5).Value = "  +0.20%  "

$ws.Cells.Item(42, 5).Value = "  -2.96%  "

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "21.57"
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -4.41%  "

$ws.Cells.Item(44, 5).Value = "  -2.84%  "

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "137.35"
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -1.04%  "

$ws.Cells.Item(47, 5).Value = "  -3.44%  "

$ws.Cells.Item(48, 5).Value = "  -3.90%  "

$ws.Cells.Item(50, 5).Value = "  -0.04%  "

$ws.Cells.Item(51, 5).Value = "  +0.04%  "
